$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the greeting text in E8 ("Good Morning" -> "GIT UPDATE")
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new selection left behind by the edit
$ws.Range("E8").Select()
